$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to place in the rightmost new diagonal column (rows 2-14, column K)
$newValues = @{
    2 = -0.2998068482859306
    3 = -0.1226260243844251
    4 = -0.7091962372939675
    5 = 0.4742753336424089
    6 = -0.4147784639390545
    7 = -1.253305011333664
    8 = 0.6169758060237021
    9 = 0.1142546880022587
    10 = 0.3970518656191074
    11 = 0.2762203743837313
    12 = -0.4272707339946972
    13 = 0.3617928071605474
    14 = 1.11229800409388
}

# Shift the staircase data block one column to the left: for each row, read
# the old values starting at column C through the row's last used column,
# then write them back starting at column B. Finally clear the vacated last
# column, and (for rows 2-14 only) populate it with the new trailing value
# that extends the staircase.
for ($r = 2; $r -le 24; $r++) {
    # Number of value-cells present in this row BEFORE the edit (rows 2-15
    # all have 10 values in columns B..K; the staircase then shrinks by one
    # column per row, from row 16 (9 values) down to row 24 (1 value)).
    if ($r -le 15) {
        $oldCount = 10
    } else {
        $oldCount = 25 - $r
    }

    $oldLastCol = 1 + $oldCount   # column A=1 plus the value columns

    # Read old values from column C through the old last column (this is one
    # fewer value than $oldCount, since column B's value is being dropped).
    $vals = @()
    for ($col = 3; $col -le $oldLastCol; $col++) {
        $vals += $ws.Cells.Item($r, $col).Value2
    }

    # Write the shifted values back starting at column B.
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value2 = $vals[$i]
    }

    # Clear the now-vacated old last column.
    $ws.Cells.Item($r, $oldLastCol).ClearContents()

    # For rows 2-14, populate that vacated column with the new trailing
    # value that extends the staircase edge.
    if ($newValues.ContainsKey($r)) {
        $ws.Cells.Item($r, $oldLastCol).Value2 = $newValues[$r]
    }
}
